$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 143, shifting the existing
# data (old rows 143-230) down to become rows 145-232.
$ws.Rows("143:144").Insert()

# --- New row 143 ---
$ws.Cells.Item(143, 1).Value = 1
$ws.Cells.Item(143, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(143, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(143, 4).Value = 44680
$ws.Cells.Item(143, 5).Value = 15
$ws.Cells.Item(143, 6).Value = "Fruta"
$ws.Cells.Item(143, 7).Value = 100102
$ws.Cells.Item(143, 8).Value = "Cítricos"
$ws.Cells.Item(143, 9).Value = 100102003
$ws.Cells.Item(143, 10).Value = "Limón"
$ws.Cells.Item(143, 11).Value = "Sutil De Gase"
$ws.Cells.Item(143, 12).Value = "Primera"
$ws.Cells.Item(143, 13).Value = 250
$ws.Cells.Item(143, 14).Value = 28000
$ws.Cells.Item(143, 15).Value = 29000
$ws.Cells.Item(143, 16).Value = 28500
$ws.Cells.Item(143, 17).Value = "$/caja 24 kilos"
$ws.Cells.Item(143, 18).Value = "Perú"
$ws.Cells.Item(143, 19).Value = 1188
$ws.Cells.Item(143, 20).Value = 24

# --- New row 144 ---
$ws.Cells.Item(144, 1).Value = 1
$ws.Cells.Item(144, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(144, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(144, 4).Value = 44680
$ws.Cells.Item(144, 5).Value = 15
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100102
$ws.Cells.Item(144, 8).Value = "Cítricos"
$ws.Cells.Item(144, 9).Value = 100102003
$ws.Cells.Item(144, 10).Value = "Limón"
$ws.Cells.Item(144, 11).Value = "Tahití"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 200
$ws.Cells.Item(144, 14).Value = 13000
$ws.Cells.Item(144, 15).Value = 14000
$ws.Cells.Item(144, 16).Value = 13500
$ws.Cells.Item(144, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(144, 18).Value = "Brasil"
$ws.Cells.Item(144, 19).Value = 750
$ws.Cells.Item(144, 20).Value = 18
